# DLAD-PART-43 edit script
# Applies the TOC-link / run-splitting / bookmark-relocation changes
# described by the target diff, using Word COM-interop primitives
# (Find, Range.Delete, Range.InsertXML).

function New-PkgXml {
    param([string]$BodyInnerXml)
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' + $BodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Returns a Range for the Nth (1-based) occurrence of $SearchText in the
# document body.
function Find-Nth {
    param($Doc, [string]$SearchText, [int]$N)
    $rng = $Doc.Content
    $rng.Start = 0
    $count = 0
    while ($true) {
        $found = $rng.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $count = $count + 1
        if ($count -eq $N) {
            return $Doc.Range($rng.Start, $rng.End)
        }
        $rng.Start = $rng.End
        $rng.End = $Doc.Content.End
    }
    throw ("Find-Nth: not found -> " + $SearchText + " occurrence " + $N)
}

# Deletes $Range's content and inserts $InnerXml (one or more <w:p>/<w:r>
# elements) in its place via InsertXML.
function Replace-RangeWithXml {
    param($Doc, $Range, [string]$InnerXml)
    $start = $Range.Start
    $Range.Delete()
    $ins = $Doc.Range($start, $start)
    $ins.InsertXML((New-PkgXml $InnerXml))
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Both "(Revised August 10, 2018 through PROCLTR 2018-16)" paragraphs:
#    split the single run into three runs, wrapping "2018" in
#    proofErr gramStart/gramEnd markers.
# ---------------------------------------------------------------------
$revisedInner = '<w:p><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">(Revised August 10, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>2018</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> through PROCLTR 2018-16)</w:t></w:r></w:p>'

$r = Find-Nth $d "(Revised August 10, 2018 through PROCLTR 2018-16)" 1
Replace-RangeWithXml $d $r $revisedInner

$r = Find-Nth $d "(Revised August 10, 2018 through PROCLTR 2018-16)" 1
Replace-RangeWithXml $d $r $revisedInner

Write-Output "step1 done"

# ---------------------------------------------------------------------
# 2) TOC hyperlink "43.102" -> split into "43." + "1" + "02" (3 runs,
#    all with rStyle=Hyperlink), still inside the same <w:hyperlink>.
# ---------------------------------------------------------------------
$hl = $d.Hyperlinks(1)
$hlRange = $d.Range($hl.Range.Start, $hl.Range.End)
$hl102Inner = '<w:p><w:hyperlink r:id="rId11" w:anchor="P43_102" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>43.</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>02</w:t></w:r></w:hyperlink></w:p>'
Replace-RangeWithXml $d $hlRange $hl102Inner

Write-Output "step2 done"

# ---------------------------------------------------------------------
# 3) TOC hyperlink "43.103" -> split into "43." + "1" + "0" + "3" (4
#    runs, all with rStyle=Hyperlink), still inside the same
#    <w:hyperlink>.
# ---------------------------------------------------------------------
$hl = $d.Hyperlinks(2)
$hlRange = $d.Range($hl.Range.Start, $hl.Range.End)
$hl103Inner = '<w:p><w:hyperlink r:id="rId12" w:anchor="P43_103" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>43.</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>0</w:t></w:r><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>3</w:t></w:r></w:hyperlink></w:p>'
Replace-RangeWithXml $d $hlRange $hl103Inner

Write-Output "step3 done"
